# "3 spaltiges layout alternative"
#
# Two changes to slide 1:
#   1. Shrink the white background rectangle ("Rechteck 14") so it is
#      3915428 EMU tall instead of 4062386 EMU (width unchanged).
#   2. Add a new picture to the slide: a duplicate of the existing
#      "Grafik 8" picture (the one using the Trump-postgis map image),
#      offset by 12pt down/right, renamed "Grafik 33".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Find-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Name -eq $name) {
            return $candidate
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) Resize "Rechteck 14"
# ---------------------------------------------------------------------
$rect = Find-ShapeByName $s "Rechteck 14"
$rect.Height = 3915428 / 914400.0 * 72.0

# ---------------------------------------------------------------------
# 2) Duplicate "Grafik 8" to create the new "Grafik 33" picture.
#
# This host assigns newly created shape ids by filling the lowest free
# gap in the slide's id space (mirroring real PowerPoint's allocator).
# The target id for the new picture is 34; to reach it deterministically
# we first consume the lower gaps with throwaway duplicates of a plain
# shape (no picture relationships involved) and discard them again, so
# the one duplicate we keep - of "Grafik 8" - lands exactly on id 34.
# ---------------------------------------------------------------------
$scratch = @()
for ($i = 0; $i -lt 11; $i++) {
    $scratch += $rect.Duplicate()
}
foreach ($tmp in $scratch) {
    $tmp.Delete()
}

$source = Find-ShapeByName $s "Grafik 8"
$newPic = $source.Duplicate()
$newPic.Name = "Grafik 33"
